$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "PDB filename" is split into "PDB or RCSB ID" (D) + "File
# Extension" (E); "Is model" survives (with a leading space) as col F; a new
# "From RCSB" column is appended as col G.
$ws.Range("D1").Value = "PDB or RCSB ID"
$ws.Range("E1").Value = "File Extension"
$ws.Range("F1").Value = " Is model"
$ws.Range("G1").Value = "From RCSB"

# Columns A/B/C are unchanged in content (narrative id / object name /
# feature id) -- just re-asserted here for clarity/completeness.
$narrativeIds = @(50569, 50569, 50569, 50569, 57196, 57197, 57198)
$objectNames  = @("JCVI_Syn3.kbase", "JCVI_Syn3.kbase", "JCVI_Syn3.kbase", "JCVI_Syn3.kbase", "JCVI_Syn3.kbase", "JCVI_Syn3.kbase", "JCVI_Syn3.kbase")
$featureIds   = @("JCVISYN3_0001", "JCVISYN3_0002", "JCVISYN3_0003", "JCVISYN3_0004", "JCVISYN3_0004", "JCVISYN3_0004", "JCVISYN3_0004")

for ($i = 0; $i -lt $narrativeIds.Count; $i++) { $ws.Cells.Item($i + 2, 1).Value = $narrativeIds[$i] }
for ($i = 0; $i -lt $objectNames.Count; $i++)  { $ws.Cells.Item($i + 2, 2).Value = $objectNames[$i] }
for ($i = 0; $i -lt $featureIds.Count; $i++)   { $ws.Cells.Item($i + 2, 3).Value = $featureIds[$i] }

# Column D (PDB/RCSB id, extension stripped) and E (file extension), split
# out of the old combined "PDB filename" values (1fat.pdb, 1nqg.pdb,
# 1fat.cif, 5o5y.pdb, 6ift, 6ifv, 6ifw).
$ws.Cells.Item(2, 4).Value = "1fat"
$ws.Cells.Item(3, 4).Value = "1nqg"
$ws.Cells.Item(2, 5).Value = "pdb"
$ws.Cells.Item(3, 5).Value = "pdb"
$ws.Cells.Item(4, 5).Value = "cif"
$ws.Cells.Item(4, 4).Value = "1fat"
$ws.Cells.Item(5, 4).Value = "5o5y"
$ws.Cells.Item(5, 5).Value = "pdb"
$ws.Cells.Item(6, 4).Value = "6ift"
$ws.Cells.Item(6, 5).Value = "pdb"
$ws.Cells.Item(7, 4).Value = "6ifv"
$ws.Cells.Item(7, 5).Value = "pdb"
$ws.Cells.Item(8, 4).Value = "6ifw"
$ws.Cells.Item(8, 5).Value = "pdb"

# Column F: "Is model" flag, same per-row values as the old "Is model" column.
$isModel = @("yes", "yes", "no", "yes", "yes", "no", "yes")
for ($i = 0; $i -lt $isModel.Count; $i++) { $ws.Cells.Item($i + 2, 6).Value = $isModel[$i] }

# Column G: new "From RCSB" column, "yes" for every row.
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r, 7).Value = "yes" }

$ws.Range("E11").Select()
